# Add new columns I ("I0") and J ("IF") to the sheet, matching the
# existing H column formatting/style, and populate the values for
# rows 2-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers "I0" and "J0"/"IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from the existing
# H1 header cell onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows 2-30: new values for columns I and J ---
$values = @{
    2  = @(7, 7)
    3  = @(9, 9)
    4  = @(8, 8)
    5  = @(6, 6)
    6  = @(7, 7)
    7  = @(10, 10)
    8  = @(6, 6)
    9  = @(7, 7)
    10 = @(6, 6)
    11 = @(5, 6)
    12 = @(8, 8)
    13 = @(6, 6)
    14 = @(7, 8)
    15 = @(9, 9)
    16 = @(10, 11)
    17 = @(5, 7)
    18 = @(6, 7)
    19 = @(9, 9)
    20 = @(8, 8)
    21 = @(5, 5)
    22 = @(8, 8)
    23 = @(6, 6)
    24 = @(6, 6)
    25 = @(6, 6)
    26 = @(4, 5)
    27 = @(5, 5)
    28 = @(5, 5)
    29 = @(9, 9)
    30 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}
